$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at L, shifting old L:M data into M:N
$ws.Columns("L").Insert(-4161)

# --- Row 1 headers ---
$ws.Range("L1").Value = "Russia"
$ws.Range("M1").Value = "Saudi Arabia"
$ws.Range("N1").Value = "USA"

# --- Row label text updates (A5, A6) ---
$ws.Range("A5").Value = "`"Governments should actively cooperate to have all countries`nconverge in terms of GDP per capita by the end of the century`""
$ws.Range("A6").Value = "Would support a global movement to tackle CC, tax millionaires,`n and fund LICs (either petition, demonstrate, strike, or donate)"

# --- Data grid values (recalculated) ---
# Row 2
$ws.Range("B2").Value = 0.556007222541788
$ws.Range("C2").Value = 0.612606004275779
$ws.Range("D2").Value = 0.60685848380173
$ws.Range("E2").Value = 0.610908673368053
$ws.Range("F2").Value = 0.723199355989802
$ws.Range("G2").Value = 0.466615113907624
$ws.Range("H2").Value = 0.582799012683891
$ws.Range("I2").Value = 0.531219428146678
$ws.Range("J2").Value = 0.531962707015654
$ws.Range("K2").Value = 0.415138078795314
$ws.Range("L2").Value = 0.595756191953926
$ws.Range("M2").Value = 0.677827112481047
$ws.Range("N2").Value = 0.48531714887618

# Row 3
$ws.Range("B3").Value = 0.50312291438834
$ws.Range("C3").Value = 0.55532380171867
$ws.Range("D3").Value = 0.561948424114003
$ws.Range("E3").Value = 0.519922424594706
$ws.Range("F3").Value = 0.539237089697176
$ws.Range("G3").Value = 0.521300773825675
$ws.Range("H3").Value = 0.566096500614007
$ws.Range("I3").Value = 0.53429638847586
$ws.Range("J3").Value = 0.380404592337776
$ws.Range("K3").Value = 0.348498257657963
$ws.Range("L3").Value = 0.604399737467109
$ws.Range("M3").Value = 0.671355992475847
$ws.Range("N3").Value = 0.431006152599421

# Row 4
$ws.Range("B4").Value = 0.680881448179833
$ws.Range("C4").Value = 0.701539116816613
$ws.Range("D4").Value = 0.710751335841921
$ws.Range("E4").Value = 0.696003924106409
$ws.Range("F4").Value = 0.737965271621624
$ws.Range("G4").Value = 0.546666497937564
$ws.Range("H4").Value = 0.721638347547208
$ws.Range("I4").Value = 0.666574918268333
$ws.Range("J4").Value = 0.652980895956235
$ws.Range("K4").Value = 0.754892621650637
$ws.Range("L4").Value = 0.688615273248795
$ws.Range("M4").Value = 0.713280127381035
$ws.Range("N4").Value = 0.608689412107398

# Row 5
$ws.Range("B5").Value = 0.609601586795904
$ws.Range("C5").Value = 0.667387097439935
$ws.Range("D5").Value = 0.580532036636799
$ws.Range("E5").Value = 0.63875633973718
$ws.Range("F5").Value = 0.749059397935025
$ws.Range("G5").Value = 0.731466745350337
$ws.Range("H5").Value = 0.710128809790465
$ws.Range("I5").Value = 0.538916979845312
$ws.Range("J5").Value = 0.558155965173253
$ws.Range("K5").Value = 0.554455577170065
$ws.Range("L5").Value = 0.778963825426238
$ws.Range("M5").Value = 0.756844881931732
$ws.Range("N5").Value = 0.448128989999106

# Row 6
$ws.Range("B6").Value = 0.675595447215337
$ws.Range("C6").Value = 0.719216740354837
$ws.Range("D6").Value = 0.68447503664547
$ws.Range("E6").Value = 0.679242709384061
$ws.Range("F6").Value = 0.792515487565139
$ws.Range("G6").Value = 0.677577969433519
$ws.Range("H6").Value = 0.746395684069502
$ws.Range("I6").Value = 0.670549607279346
$ws.Range("J6").Value = 0.648965346480897
$ws.Range("K6").Value = 0.531016915460678
$ws.Range("L6").ClearContents()
$ws.Range("M6").Value = 0.727098526374066
$ws.Range("N6").Value = 0.646544698767764

# Row 7
$ws.Range("B7").Value = 0.364717906507653
$ws.Range("C7").Value = 0.40721776360235
$ws.Range("D7").Value = 0.41037218564678
$ws.Range("E7").Value = 0.376640297322738
$ws.Range("F7").Value = 0.438569657841468
$ws.Range("G7").Value = 0.245046285790811
$ws.Range("H7").Value = 0.412435935523674
$ws.Range("I7").Value = 0.388292136341869
$ws.Range("J7").Value = 0.320423008221701
$ws.Range("K7").Value = 0.208534721957305
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = 0.336970942498412

# Row 8
$ws.Range("B8").Value = 0.347853243460036
$ws.Range("C8").Value = 0.389632627264691
$ws.Range("D8").Value = 0.344637035841507
$ws.Range("E8").Value = 0.330038030870184
$ws.Range("F8").Value = 0.541474107987616
$ws.Range("G8").ClearContents()
$ws.Range("H8").Value = 0.391565913028478
$ws.Range("I8").Value = 0.357922697953927
$ws.Range("J8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = 0.285582964276425

# Row 9
$ws.Range("B9").Value = 0.410626908494325
$ws.Range("C9").Value = 0.439007996245569
$ws.Range("D9").Value = 0.313980111032844
$ws.Range("E9").Value = 0.439081145937391
$ws.Range("F9").Value = 0.476354840470821
$ws.Range("G9").Value = 0.378195194728228
$ws.Range("H9").Value = 0.472874811326255
$ws.Range("I9").Value = 0.440260307177303
$ws.Range("J9").Value = 0.398448398482951
$ws.Range("K9").Value = 0.309849031423241
$ws.Range("L9").Value = 0.389233362357354
$ws.Range("M9").Value = 0.662536527281334
$ws.Range("N9").Value = 0.378969990866226

